$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.Value = "'28.826.60"
$cD.Style = "Normal"
$cE = $ws.Range("E2")
$cE.Value = "  +0.66%  "

$cD = $ws.Range("D3")
$cD.Value = "'1.891.29"
$cD.Style = "Normal"
$cE = $ws.Range("E3")
$cE.Value = "  +1.32%  "

$cD = $ws.Range("D4")
$cD.Value = "'1.004"
$cD.Style = "Normal"
$cE = $ws.Range("E4")
$cE.Value = "  -0.23%  "

$cD = $ws.Range("D5")
$cD.Value = "'326.12"
$cD.Style = "Normal"
$cE = $ws.Range("E5")
$cE.Value = "  -0.13%  "

$cD = $ws.Range("D6")
$cD.Value = "'1.004"
$cD.Style = "Normal"
$cE = $ws.Range("E6")
$cE.Value = "  -0.08%  "

$cD = $ws.Range("D7")
$cD.Value = "'0.4575"
$cD.Style = "Normal"
$cE = $ws.Range("E7")
$cE.Value = "  -1.24%  "

$cD = $ws.Range("D8")
$cD.Value = "'0.3856"
$cD.Style = "Normal"
$cE = $ws.Range("E8")
$cE.Value = "  -1.33%  "

$cD = $ws.Range("D9")
$cD.Value = "'0.07856"
$cD.Style = "Normal"
$cE = $ws.Range("E9")
$cE.Value = "  -0.83%  "

$cD = $ws.Range("D10")
$cD.Value = "'0.9966"
$cD.Style = "Normal"
$cE = $ws.Range("E10")
$cE.Value = "  +2.85%  "

$cD = $ws.Range("D11")
$cD.Value = "'21.68"
$cD.Style = "Normal"
$cE = $ws.Range("E11")
$cE.Value = "  -2.77%  "

$cD = $ws.Range("D12")
$cD.Value = "'1.888.91"
$cD.Style = "Normal"
$cE = $ws.Range("E12")
$cE.Value = "  +1.43%  "

$cD = $ws.Range("D13")
$cD.Value = "'6.978"
$cD.Style = "Normal"
$cE = $ws.Range("E13")
$cE.Value = "  +0.64%  "

$cD = $ws.Range("D14")
$cD.Value = "'5.682"
$cD.Style = "Normal"
$cE = $ws.Range("E14")
$cE.Value = "  -0.81%  "

$cD = $ws.Range("D15")
$cD.Value = "'0.06942"
$cD.Style = "Normal"
$cE = $ws.Range("E15")
$cE.Value = "  -0.29%  "

$cD = $ws.Range("D16")
$cD.Value = "'87.10"
$cD.Style = "Normal"
$cE = $ws.Range("E16")
$cE.Value = "  -1.24%  "

$cD = $ws.Range("D17")
$cD.Value = "'1.004"
$cD.Style = "Normal"
$cE = $ws.Range("E17")
$cE.Value = "  -0.20%  "

$cD = $ws.Range("D18")
$cD.Value = "'0.00001002"
$cD.Style = "Normal"
$cE = $ws.Range("E18")
$cE.Value = "  -0.44%  "

$cD = $ws.Range("D19")
$cD.Value = "'16.85"
$cD.Style = "Normal"
$cE = $ws.Range("E19")
$cE.Value = "  -0.68%  "

$cD = $ws.Range("D20")
$cD.Value = "'1.007"
$cD.Style = "Normal"
$cE = $ws.Range("E20")
$cE.Value = "  +0.22%  "

$cD = $ws.Range("D21")
$cD.Value = "'28.839.29"
$cD.Style = "Normal"
$cE = $ws.Range("E21")
$cE.Value = "  +0.71%  "

$cD = $ws.Range("D22")
$cD.Value = "'5.285"
$cD.Style = "Normal"
$cE = $ws.Range("E22")
$cE.Value = "  -0.70%  "

$cD = $ws.Range("D23")
$cD.Value = "'10.95"
$cD.Style = "Normal"
$cE = $ws.Range("E23")
$cE.Value = "  -1.17%  "

$cD = $ws.Range("D24")
$cD.Value = "'2.147.29"
$cD.Style = "Normal"
$cE = $ws.Range("E24")
$cE.Value = "  +2.57%  "

$cD = $ws.Range("D25")
$cD.Value = "'2.072"
$cD.Style = "Normal"
$cE = $ws.Range("E25")
$cE.Value = "  -2.23%  "

$cD = $ws.Range("D26")
$cD.Value = "'154.45"
$cD.Style = "Normal"
$cE = $ws.Range("E26")
$cE.Value = "  +0.50%  "

$cD = $ws.Range("D27")
$cD.Value = "'19.22"
$cD.Style = "Normal"
$cE = $ws.Range("E27")
$cE.Value = "  -0.55%  "

$cD = $ws.Range("D28")
$cD.Value = "'5.741"
$cD.Style = "Normal"
$cE = $ws.Range("E28")
$cE.Value = "  +0.73%  "

$cD = $ws.Range("D29")
$cD.Value = "'118.05"
$cD.Style = "Normal"
$cE = $ws.Range("E29")
$cE.Value = "  -1.11%  "

$cD = $ws.Range("D30")
$cD.Value = "'1.908"
$cD.Style = "Normal"
$cE = $ws.Range("E30")
$cE.Value = "  -4.54%  "

$cD = $ws.Range("D31")
$cD.Value = "'0.09306"
$cD.Style = "Normal"
$cE = $ws.Range("E31")
$cE.Value = "  -0.71%  "

$cD = $ws.Range("D32")
$cD.Value = "'0.9128"
$cD.Style = "Normal"
$cE = $ws.Range("E32")
$cE.Value = "  -1.94%  "

$cD = $ws.Range("D33")
$cD.Value = "'5.298"
$cD.Style = "Normal"
$cE = $ws.Range("E33")
$cE.Value = "  -0.53%  "

$cD = $ws.Range("D34")
$cD.Value = "'1.329"
$cD.Style = "Normal"
$cE = $ws.Range("E34")
$cE.Value = "  -1.14%  "

$cD = $ws.Range("D35")
$cD.Value = "'3.247"
$cD.Style = "Normal"
$cE = $ws.Range("E35")
$cE.Value = "  -3.25%  "

$cD = $ws.Range("D36")
$cD.Value = "'0.05690"
$cD.Style = "Normal"
$cE = $ws.Range("E36")
$cE.Value = "  -2.53%  "

$cD = $ws.Range("D37")
$cD.Value = "'1.155"
$cD.Style = "Normal"
$cE = $ws.Range("E37")
$cE.Value = "  +0.43%  "

$cD = $ws.Range("D38")
$cD.Value = "'0.02050"
$cD.Style = "Normal"
$cE = $ws.Range("E38")
$cE.Value = "  -3.60%  "

$cD = $ws.Range("D39")
$cD.Value = "'7.672"
$cD.Style = "Normal"
$cE = $ws.Range("E39")
$cE.Value = "  -2.97%  "

$cD = $ws.Range("D40")
$cD.Value = "'0.5594"
$cD.Style = "Normal"
$cE = $ws.Range("E40")
$cE.Value = "  -1.09%  "

$cD = $ws.Range("D41")
$cD.Value = "'0.1777"
$cD.Style = "Normal"
$cE = $ws.Range("E41")
$cE.Value = "  -0.32%  "

$cD = $ws.Range("D42")
$cD.Value = "'9.661"
$cD.Style = "Normal"
$cE = $ws.Range("E42")
$cE.Value = "  -2.64%  "

$cD = $ws.Range("D43")
$cD.Value = "'0.07174"
$cD.Style = "Normal"
$cE = $ws.Range("E43")
$cE.Value = "  -0.93%  "

$cD = $ws.Range("D44")
$cD.Value = "'0.5281"
$cD.Style = "Normal"
$cE = $ws.Range("E44")
$cE.Value = "  -0.68%  "

$cD = $ws.Range("D45")
$cD.Value = "'11.48"
$cD.Style = "Normal"
$cE = $ws.Range("E45")
$cE.Value = "  -2.12%  "

$cD = $ws.Range("D46")
$cD.Value = "'2.136"
$cD.Style = "Normal"
$cE = $ws.Range("E46")
$cE.Value = "  -1.17%  "

$cD = $ws.Range("D47")
$cD.Value = "'1.113"
$cD.Style = "Normal"
$cE = $ws.Range("E47")
$cE.Value = "  -1.95%  "

$cD = $ws.Range("D48")
$cD.Value = "'1.814"
$cD.Style = "Normal"
$cE = $ws.Range("E48")
$cE.Value = "  -1.69%  "

$cD = $ws.Range("D49")
$cD.Value = "'112.05"
$cD.Style = "Normal"
$cE = $ws.Range("E49")
$cE.Value = "  -1.26%  "

$cD = $ws.Range("D50")
$cD.Value = "'2.455"
$cD.Style = "Normal"
$cE = $ws.Range("E50")
$cE.Value = "  +4.69%  "

$cD = $ws.Range("D51")
$cD.Value = "'1.003"
$cD.Style = "Normal"
$cE = $ws.Range("E51")
$cE.Value = "  -0.05%  "

